# "Update Shadow Rate with Latest Data"
# Refreshes the fedfundsrate_shadow (shadow rate) series in column C for the
# existing effective-lower-bound period (rows 81-128 and 147-151) with
# re-estimated model output, and appends one new observation (2021.5) as
# row 152.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revised shadow-rate estimates for existing rows -----------------------
$updates = @{
    81  = 0.99639182326900499
    82  = 1.0029176028518938
    83  = 1.0095228599094508
    84  = 1.4328447041452641
    85  = 1.9495340410316864
    86  = 2.4695795576779034
    87  = 2.9429727164916653
    88  = 3.4597066708629276
    89  = 3.9797765655578132
    90  = 4.4565120716320239
    91  = 4.9065770706629763
    92  = 5.2466363057630083
    93  = 5.2466883505500217
    94  = 5.2567322996835308
    95  = 5.2501011760757832
    96  = 5.0734614948653389
    97  = 4.4968131695784175
    98  = 3.176822794968337
    99  = 2.0868253592999597
    100 = 1.9401557652927393
    101 = 0.50775031191936382
    102 = 1.6937772700994946
    103 = 0.67587175210195305
    104 = -0.10559049935601195
    105 = -0.17776334042235442
    106 = -0.1710687472153638
    107 = -1.2615725352793983
    108 = -1.138310679947796
    109 = -1.5852699405839221
    110 = -1.3650766135385073
    111 = -1.0773344011158592
    112 = -1.8490520193814364
    113 = -1.2301890236138813
    114 = -1.9960226293425531
    115 = -1.7099618130806626
    116 = -1.3547201128070174
    117 = -2.455059325937381
    118 = -1.1911912825350535
    119 = -0.70237546625840919
    120 = -0.47871997828423662
    121 = -0.56565540430401073
    122 = -0.77565234754822798
    123 = -0.71476232806565454
    124 = -0.46028858753257573
    125 = -0.21487205080273908
    126 = 0.17503399903986416
    127 = 0.035722522054393302
    128 = 0.053568438689199915
    147 = 7.8030407161213233
    148 = -4.7399670312971214
    149 = -3.1478361323883819
    150 = -3.0435707193201766
    151 = -2.7772534156314532
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# --- Append the newest observation (2021Q3 / 2021.5) -----------------------
$ws.Cells.Item(152, 1).Value = 2021.5
$ws.Cells.Item(152, 2).Value = 0
$ws.Cells.Item(152, 3).Value = -1.7280701584331148
